# Fruta / hortaliza, semanal
# Prepend a new weekly price-update block (3 rows) for Chirimoya, pushing the
# existing data down by 3 rows (the sheet has no formulas referencing these
# rows, so a plain row-insert + fill is safe).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above the current row 159 (existing rows 159:177
# shift down to 162:180).
$ws.Range("A159:A161").EntireRow.Insert()

# Shared "template" values that are constant across every row of this
# Comercializadora del Agro de Limarí / Chirimoya block.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$prodId    = 100107
$producto  = "Otros"
$catId     = 100107002
$categoria = "Chirimoya"
$variedad  = "Cultivar IV Región"
$unidad    = "`$/bandeja 10 kilos"
$origen    = "Provincia de Limarí"

# New rows: date 45204 (2023-10-05), qualities Especial / Primera / Segunda.
$rows = @(
    @{ Row = 159; Calidad = "Especial"; Volumen = 300; PMin = 20000; PMax = 21000; PProm = 20500; PKg = 2050 },
    @{ Row = 160; Calidad = "Primera";  Volumen = 400; PMin = 17000; PMax = 18000; PProm = 17500; PKg = 1750 },
    @{ Row = 161; Calidad = "Segunda";  Volumen = 360; PMin = 14000; PMax = 15000; PProm = 14500; PKg = 1450 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = 45204
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $prodId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $catId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = 10
}
